$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '43.361.57'
$ws.Range("E2").Value = '  +2.89%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '2.309.40'
$ws.Range("E3").Value = '  +4.35%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  +0.05%  '

# Row 5 - BNB
$ws.Range("D5").Value = "'252.83"
$ws.Range("E5").Value = '  +0.29%  '

# Row 6 - XRP
$ws.Range("D6").Value = "'0.642"
$ws.Range("E6").Value = '  +3.22%  '

# Row 7 - Solana
$ws.Range("D7").Value = "'74.28"
$ws.Range("E7").Value = '  +9.37%  '

# Row 8 - USDC
$ws.Range("E8").Value = '  -0.06%  '

# Row 9 - Cardano
$ws.Range("D9").Value = "'0.641"
$ws.Range("E9").Value = '  +3.67%  '

# Row 10 - Avalanche
$ws.Range("D10").Value = "'39.46"
$ws.Range("E10").Value = '  +1.28%  '

# Row 11 - Dogecoin
$ws.Range("D11").Value = "'0.0984"
$ws.Range("E11").Value = '  +4.89%  '

# Row 12 - OKB
$ws.Range("D12").Value = "'59.41"
$ws.Range("E12").Value = '  +0.00%  '

# Row 13 - Polkadot
$ws.Range("D13").Value = "'7.43"
$ws.Range("E13").Value = '  +5.77%  '

# Row 14 - TRON
$ws.Range("E14").Value = '  +1.70%  '

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '2.654.87'
$ws.Range("E15").Value = '  +4.35%  '

# Row 16 - Chainlink
$ws.Range("D16").Value = "'15.54"
$ws.Range("E16").Value = '  +7.22%  '

# Row 17 - Polygon
$ws.Range("D17").Value = "'0.881"
$ws.Range("E17").Value = '  +1.22%  '

# Row 18 - WrappedEther
$ws.Range("D18").Value = '2.309.57'
$ws.Range("E18").Value = '  +4.99%  '

# Row 19 - WrappedBTC
$ws.Range("D19").Value = '43.248.65'
$ws.Range("E19").Value = '  +2.96%  '

# Row 20 - ShibaInu
$ws.Range("E20").Value = '  +4.55%  '

# Row 21 - Uniswap
$ws.Range("E21").Value = '  +3.53%  '

# Row 22 - Litecoin
$ws.Range("D22").Value = "'72.87"
$ws.Range("E22").Value = '  +0.85%  '

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "'236.30"
$ws.Range("E23").Value = '  +2.18%  '

# Row 24 - ImmutableX
$ws.Range("E24").Value = '  +10.43%  '

# Row 25 - WEMIXToken
$ws.Range("D25").Value = "'3.91"
$ws.Range("E25").Value = '  +1.05%  '

# Row 26 - Cosmos
$ws.Range("D26").Value = "'11.64"
$ws.Range("E26").Value = '  +4.29%  '

# Row 27 - Dai
$ws.Range("E27").Value = '  -0.19%  '

# Row 28 - PancakeSwap
$ws.Range("D28").Value = "'2.43"
$ws.Range("E28").Value = '  +1.01%  '

# Row 29 - LEO
$ws.Range("D29").Value = "'3.65"
$ws.Range("E29").Value = '  -1.54%  '

# Row 30 - Toncoin
$ws.Range("E30").Value = '  -0.03%  '

# Row 31 - Monero
$ws.Range("D31").Value = "'167.77"
$ws.Range("E31").Value = '  +0.49%  '

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "'21.19"
$ws.Range("E32").Value = '  +3.85%  '

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "'6.41"
$ws.Range("E33").Value = '  +7.94%  '

# Row 34 - Kaspa
$ws.Range("E34").Value = '  +5.88%  '

# Row 35 - now Hedera (was InjectiveProtocol)
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = "'0.0815"
$ws.Range("E35").Value = '  +4.88%  '

# Row 36 - now InjectiveProtocol (was Hedera)
$ws.Range("B36").Value = 'InjectiveProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D36").Value = "'31.70"
$ws.Range("E36").Value = '  +21.86%  '

# Row 37 - Stellar
$ws.Range("E37").Value = '  +3.78%  '

# Row 38 - RenderToken
$ws.Range("E38").Value = '  +12.07%  '

# Row 39 - Filecoin
$ws.Range("D39").Value = "'4.79"
$ws.Range("E39").Value = '  +4.39%  '

# Row 40 - VeChain
$ws.Range("E40").Value = '  -0.77%  '

# Row 41 - Celestia
$ws.Range("D41").Value = "'14.46"
$ws.Range("E41").Value = '  +20.68%  '

# Row 42 - LidoDAOToken
$ws.Range("D42").Value = "'2.37"
$ws.Range("E42").Value = '  +6.52%  '

# Row 43 - THORChain
$ws.Range("D43").Value = "'6.03"
$ws.Range("E43").Value = '  +6.42%  '

# Row 44 - Algorand
$ws.Range("D44").Value = "'0.217"
$ws.Range("E44").Value = '  +11.12%  '

# Row 45 - now FraxShare (was MultiversX)
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = "'9.18"
$ws.Range("E45").Value = '  +7.49%  '

# Row 46 - now MultiversX (was FraxShare)
$ws.Range("B46").Value = 'MultiversX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D46").Value = "'62.40"
$ws.Range("E46").Value = '  +1.85%  '

# Row 47 - FTXToken
$ws.Range("D47").Value = "'4.90"
$ws.Range("E47").Value = '  -4.07%  '

# Row 48 - Cronos
$ws.Range("E48").Value = '  +3.86%  '

# Row 49 - ARBITRUM
$ws.Range("E49").Value = '  +3.53%  '

# Row 50 - BinanceUSD
$ws.Range("E50").Value = '  +0.18%  '

# Row 51 - Aave
$ws.Range("D51").Value = "'99.11"
$ws.Range("E51").Value = '  +6.85%  '
